$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.395439028739929
$ws.Range("B1").Value = 1.84431004524231
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 1.903853535652161
$ws.Range("E1").Value = 0.7334681749343872
